# Atualizacao rapida de agenda as 9:37:44,98
# Rebuilds the "agenda" rows: removes the old Roberto/Pedro job rows,
# inserts the new Giovani job list (rows 2-11), keeps Roberto's first row
# (now pointing at a different job) and pushes Pedro's remaining two rows
# down to 12-13. Clears the (now unused) "Observação"/"Cobrança"/"Status"
# columns for every data row and appends 7 fresh blank rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 : Roberto / Fazenda Mata da Pedra -------------------------------
$ws.Range("A2").Value = "Roberto"
$ws.Range("B2").Value = "'0629"
$ws.Range("C2").Value = "Fazenda Mata da Pedra"
$ws.Range("D2").Value = "Uma câmera sem funcionar"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# --- rows 3-11 : Giovani's job list -----------------------------------------
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0288"
$ws.Range("C3").Value = "Viaceu Loja"
$ws.Range("D3").Value = "Sem comunicação de câmeras."
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0304"
$ws.Range("C4").Value = "Cimentão"
$ws.Range("D4").Value = "Trocar meio para DDNS (era antes)."
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0867"
$ws.Range("C5").Value = "Rotoplast"
$ws.Range("D5").Value = "Todas as câmeras com a logo do fabricante."
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""

$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0706"
$ws.Range("C6").Value = "Lar das Meninnas"
$ws.Range("D6").Value = "4/5 câmeras com a logo do fabricante."
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""

$ws.Range("A7").Value = "Giovani"
$ws.Range("B7").Value = "'0885"
$ws.Range("C7").Value = "Arcelormital ADM"
$ws.Range("D7").Value = "Parece que tem um sensor caído, ninja mandou foto pra nós."
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""

$ws.Range("A8").Value = "Giovani"
$ws.Range("B8").Value = "'0887"
$ws.Range("C8").Value = "Arcelomital Galpão"
$ws.Range("D8").Value = "Zona 3 parece que tá aberta."
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""

$ws.Range("A9").Value = "Giovani"
$ws.Range("B9").Value = "'0868"
$ws.Range("C9").Value = "Tricostura"
$ws.Range("D9").Value = "Reposiconar duas câmeras."
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""

$ws.Range("A10").Value = "Giovani"
$ws.Range("B10").Value = "'0351"
$ws.Range("C10").Value = "Oribes Batista"
$ws.Range("D10").Value = "Instalação de um DVR novo."
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""

$ws.Range("A11").Value = "Giovani"
$ws.Range("B11").Value = "'0580"
$ws.Range("C11").Value = "Toyota Oficina"
$ws.Range("D11").Value = "Central não reporta armado nem no buffer."
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""

# --- rows 12-13 : Pedro (itaúna)'s remaining jobs, shifted down -------------
$ws.Range("A12").Value = "Pedro (itaúna)"
$ws.Range("B12").Value = "2138"
$ws.Range("C12").Value = "Guia CWK"
$ws.Range("D12").Value = "Local sem comunicação de alarmes."
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""

$ws.Range("A13").Value = "Pedro (itaúna)"
$ws.Range("B13").Value = "2565"
$ws.Range("C13").Value = "Rodonaves"
$ws.Range("D13").Value = "Reparo e acesso em câmeras não monitoradas."
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""

# row 2 no longer wraps a long "Observação" line, so its custom height goes
# back to the sheet's default.
$ws.Rows(2).AutoFit()

# --- append 7 blank rows at the end of the table (145-151) -----------------
# Copy the formatting of the last existing blank row so the new rows reuse
# the same cell styles/borders instead of minting new ones.
$lastRow = 144
$newFirst = $lastRow + 1
$newLast = $lastRow + 7
$ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
$ws.Range("A" + $newFirst + ":G" + $newLast).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- refresh the view so the active cell matches where the author ended up -
$ws.Range("G12").Select()
